# Page Speed insight - FCP & LCP
# Clear the Desktop LCP/FCP sample values out of row 2 (E2:F2) and narrow
# those two columns down, then leave the selection on D2 (matches the
# author's commit which trims the desktop FCP/LCP sample for the first row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the E2 / F2 values (shared-string cells "3 s" / "1.5 s").
$ws.Range("E2:F2").ClearContents()

# Shrink columns E and F - they no longer need to fit the removed values.
$ws.Columns.Item(5).ColumnWidth = 3.83
$ws.Columns.Item(6).ColumnWidth = 3.83

# Move the active selection to D2.
$ws.Range("D2").Select()
